# Edit script: Add "Shale Talents Reworked" (STR) talent strings to the workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from Sheet2 to Sheet1 (defined name strings_1 auto-updates its reference)
$ws.Name = "Sheet1"

# Add the new talent id/text rows (106-138), copying formatting from the
# last existing "talent description" style row so colors / wrap-text / number format match
$ws.Range("A102:B102").Copy()
$ws.Range("A106:B138").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A106").Value = 6610104
$ws.Range("B106").Value = "Shale activates a mode that focuses on offense, taking penalties to defense and armor in exchange for bonuses to damage and attack. Slam, Quake, and Killing Blow add bonuses to this mode."

$ws.Range("A107").Value = 6610105
$ws.Range("B107").Value = "Shale slams a stony fist into the enemy target. If the blow connects, it generates an automatic critical hit and knocks the target back. The target is also stunned briefly unless it passes a physical resistance check. After learning this talent, Shale gains bonuses to damage, attack, and armor penetration whenever Pulverizing Blows is active."
$ws.Rows("107").RowHeight = 22.5

$ws.Range("A108").Value = 6610106
$ws.Range("B108").Value = "Shale strikes the ground repeatedly, sending out multiple shockwaves. Nearby creatures take damage and are knocked to the ground unless they pass a physical resistance check. They are also slowed briefly unless they pass another physical resistance check. After learning this talent, Shale gains a bonus to movement speed, while suffering a greater penalty to defense, whenever Pulverizing Blows is active."
$ws.Rows("108").RowHeight = 22.5

$ws.Range("A109").Value = 6610107
$ws.Range("B109").Value = "Shale focuses strength into one devastating blow against an enemy. This attack gains bonus armor penetration and is an automatic critical hit if it connects. If the target has less than 50% health, this attack inflicts increasing damage as the target’s health diminishes. Shale suffers penalties to attack and stamina regeneration for a short time after using this attack. After learning this talent, Shale gains bonuses to damage, attack, armor penetration, and critical chance whenever Pulverizing Blows is active."
$ws.Rows("109").RowHeight = 33.75

$ws.Range("A110").Value = 6610108
$ws.Range("B110").Value = "Shale activates a mode that focuses on defense, gaining bonuses to armor, physical, and elemental resistances, while suffering penalties to damage and critical chance. Whenever this mode is active, Shale’s threat toward all nearby enemies increases over time, making them more likely to target Shale. Bellow, Stone Will, and Regenerating Burst add bonuses to this mode."
$ws.Rows("110").RowHeight = 22.5

$ws.Range("A111").Value = 6610109
$ws.Range("B111").Value = "Shale erupts with the sound of the Rock itself, stunning nearby enemies unless they pass a mental resistance check against Shale’s constitution. They are also dazed, and suffer penalties to attack and defense, unless they pass another mental resistance check against Shale’s constitution. Enemies stunned by this attack cannot resist being dazed. The power of this talent depends on Shale's constitution. This talent increases the bonuses of Stoneheart."
$ws.Rows("111").RowHeight = 33.75

$ws.Range("A112").Value = 6610110
$ws.Range("B112").Value = "Shale attracts the attention of a single enemy, who now views Shale as the most urgent threat on the battlefield and is likely to target Shale. After this talent is learned, Stoneheart will also give Shale bonuses to health regeneration and stamina regeneration, and enemies are even more likely to target Shale when that mode is active."
$ws.Rows("112").RowHeight = 22.5

$ws.Range("A113").Value = 6610111
$ws.Range("B113").Value = "Shale explodes with energy, damaging all nearby enemies. Enemies hit by the explosion are knocked down or knocked back unless they pass a physical resistance check against Shale’s constitution. They are also stunned unless they pass a mental resistance check against Shale’s constitution. The energy gives Shale a burst of health regeneration and stamina regeneration. The lower Shale's health or stamina, the stronger the regeneration. The power of this talent also depends on Shale's constitution. This talent further increases the bonuses of Stoneheart."
$ws.Rows("113").RowHeight = 33.75

$ws.Range("A114").Value = 6610112
$ws.Range("B114").Value = "Shale activates a mode that enables long-range attacks, gaining a large defense bonus against missile attacks while suffering penalties to attack and movement speed. While in this mode, Shale emanates an aura that grants nearby party members bonuses to ranged critical chance and ranged attack speed. Hurl Rock, Earthen Grasp, and Rock Barrage add bonuses to this mode."
$ws.Rows("114").RowHeight = 22.5

$ws.Range("A115").Value = 6610113
$ws.Range("B115").Value = "Shale pulls a rock from the ground and hurls it to a location, dealing physical damage to all creatures in the impact area. Creatures closer to the center of the impact take more damage, and those within 3 meters of the center are also knocked down unless they pass a physical resistance check against Shale's strength. Friendly fire possible. This talent increases the bonuses for party members and Shale’s missile deflection."
$ws.Rows("115").RowHeight = 22.5

$ws.Range("A116").Value = 6610114
$ws.Range("B116").Value = "Shale pounds the earth, immobilizing enemies unless they pass a physical resistance check, in which case they suffer penalties to attack speed and movement speed instead. Shale’s willpower increases the duration of the effects. This talent further increases Shale’s missile deflection."
$ws.Rows("116").RowHeight = 22.5

$ws.Range("A117").Value = 6610115
$ws.Range("B117").Value = "Shale tosses up multiple rocks that crash down in the targeted area. Creatures within the area take damage and are knocked down or knocked back unless they pass a physical resistance check. Friendly fire possible. This talent increases the bonuses for party members and Shale’s missile deflection. Nearby party members now also gain bonuses to attack and damage when using ranged weapons. Additionally, Shale’s aura can now shield others, granting party members within 3 meters a bonus to missile deflection."
$ws.Rows("117").RowHeight = 33.75

$ws.Range("A118").Value = 6610116
$ws.Range("B118").Value = "Shale activates a support mode that imbues nearby party members with bonuses to defense, armor, and all resistances. Shale is immobilized when in this mode, suffering a penalty to defense while gaining bonuses to armor, physical, spell, and elemental resistances. Inner Reserves, Renewed Assault, and Supernatural Resilience add bonuses to this mode."
$ws.Rows("118").RowHeight = 22.5

$ws.Range("A119").Value = 6610117
$ws.Range("B119").Value = "Whenever Stone Aura is active, party members within the aura now receive bonuses to health regeneration, stamina regeneration, and spellpower. This talent also increases the radius of Stone Aura."

$ws.Range("A120").Value = 6610118
$ws.Range("B120").Value = "Whenever Stone Aura is active, Shale gains additional bonuses to armor, spell, and elemental resistances. Party members within the aura gain additional bonuses to health regeneration, stamina regeneration, and spellpower, as well as a bonus to movement speed. This talent also increases the radius of Stone Aura."
$ws.Rows("120").RowHeight = 22.5

$ws.Range("A121").Value = 6610119
$ws.Range("B121").Value = "Whenever Stone Aura is active, Shale gains additional bonuses to armor, physical, spell, and elemental resistances. Party members within the aura gain additional bonuses to defense, armor, and all resistances. The aura now also grants additional bonuses to party members when they stay close to Shale. This talent also increases the radius of Stone Aura."
$ws.Rows("121").RowHeight = 22.5

$ws.Range("A122").Value = 6610120
$ws.Range("B122").Value = "Bonuses to damage and attack; penalties to defense and armor (Slam: bonus to armor penetration; Quake: bonus to movement speed)."

$ws.Range("A123").Value = 6610121
$ws.Range("B123").Value = "Knocked down; possible penalty to movement speed."

$ws.Range("A124").Value = 6610122
$ws.Range("B124").Value = "Shale suffers reduced attack and stamina regeneration. Target is knocked down or knocked back."

$ws.Range("A125").Value = 6610123
$ws.Range("B125").Value = "Bonuses to armor, physical, and elemental resistances; penalties to damage and critical chance (Stone Will: bonus to health regeneration and stamina regeneration)."

$ws.Range("A126").Value = 6610124
$ws.Range("B126").Value = "Stunned; possible penalties to attack and defense."

$ws.Range("A127").Value = 6610125
$ws.Range("B127").Value = "Shale gains increased health regeneration and stamina regeneration. Enemies are knocked down, knocked back, or stunned."

$ws.Range("A128").Value = 6610126
$ws.Range("B128").Value = "Shale gains increased missile deflection, but suffers reduced movement speed and attack. Party members gain bonuses to ranged critical chance and aim speed."

$ws.Range("A129").Value = 6610127
$ws.Range("B129").Value = "Paralyzed; penalties to attack speed and movement speed."

$ws.Range("A130").Value = 6610128
$ws.Range("B130").Value = "Knocked down or knocked back."

$ws.Range("A131").Value = 6610129
$ws.Range("B131").Value = "Shale is paralyzed but gains bonuses to armor and resistances. Party members gain bonuses to defense, armor, and resistances."

$ws.Range("A132").Value = 6610130
$ws.Range("B132").Value = "Bonuses to ranged attack and damage; possible bonus to missile deflection."

$ws.Range("A133").Value = 6610131
$ws.Range("B133").Value = "Bonuses to defense, armor, resistances, and spellpower."

$ws.Range("A134").Value = 6610132
$ws.Range("B134").Value = "Stone Will"

$ws.Range("A135").Value = 6610133
$ws.Range("B135").Value = "With a will of stone, Shale becomes a living fortress. All active movement speed modifiers are immediately removed. For a time, Shale gains increased physical and mental resistances, and becomes immune to knockdown, knockback, and slip effects. All incoming damage is partially absorbed, until a maximum amount is reached. Shale is slowed during this time, but will also shrug off any incoming effect that would alter movement speed. The power of this talent depends on Shale's constitution. After this talent is learned, Shale gains bonuses to health regeneration and stamina regeneration whenever Stoneheart is active, and enemies are even more likely to target Shale."
$ws.Rows("135").RowHeight = 45

$ws.Range("A136").Value = 6610134
$ws.Range("B136").Value = "Bonuses to physical and mental resistances; incoming damage is partially absorbed; immune to knockdown, knockback, slip, and movement speed modifications."

$ws.Range("A137").Value = 6610135
$ws.Range("B137").Value = "With a will of stone, Shale becomes a living fortress. All active movement speed modifiers are immediately removed. For a time, Shale gains increased physical and mental resistances, and becomes immune to knockdown, stun, and slip effects. All incoming damage is partially absorbed, until a maximum amount is reached. Shale is slowed during this time, but will also shrug off any incoming effect that would alter movement speed. The power of this talent depends on Shale's constitution. After this talent is learned, Shale gains bonuses to health regeneration and stamina regeneration whenever Stoneheart is active, and enemies are even more likely to target Shale."
$ws.Rows("137").RowHeight = 45

$ws.Range("A138").Value = 6610136
$ws.Range("B138").Value = "Bonuses to physical and mental resistances; incoming damage is partially absorbed; immune to knockdown, stun, slip, and movement speed modifications."

# New comment on A106 marking the start of the Shale Talents Reworked (STR) block
$ws.Range("A106").AddComment("Shale Talents Reworked")

# Column B got wider to accommodate the new, longer talent descriptions
$ws.Columns("B").ColumnWidth = 152

# Move the active selection the way it was left in the authored workbook
$ws.Range("B148").Select()

Write-Output "Shale Talents Reworked strings added"
